# Kuwait GDP per Capita (CCode 414): refresh the series with the latest
# Clio-Infra release and extend the table through 2016.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column E stores the GDP-per-capita figures as text (the workbook keeps
# them as strings even though they look numeric), so force a Text format
# before writing the values to stop Excel from re-typing them as numbers.
$ws.Range("E2:E68").NumberFormat = "@"

# row => updated GDP per Capita value (years 1950-2010 already on the sheet)
$gdpUpdates = @{
    2  = "46031"  # 1950
    3  = "47464"  # 1951
    4  = "47856"  # 1952
    5  = "49989"  # 1953
    6  = "52920"  # 1954
    7  = "51417"  # 1955
    8  = "52404"  # 1956
    9  = "50126"  # 1957
    10  = "47671"  # 1958
    11  = "47131"  # 1959
    12  = "45929"  # 1960
    13  = "41622"  # 1961
    14  = "42150"  # 1962
    15  = "40377"  # 1963
    16  = "40333"  # 1964
    17  = "37511"  # 1965
    18  = "38335"  # 1966
    19  = "35720"  # 1967
    20  = "35546"  # 1968
    21  = "33415"  # 1969
    22  = "48927"  # 1970
    23  = "49302"  # 1971
    24  = "48283"  # 1972
    25  = "42542"  # 1973
    26  = "34962"  # 1974
    27  = "28950"  # 1975
    28  = "28956"  # 1976
    29  = "26168"  # 1977
    30  = "26363"  # 1978
    31  = "28174"  # 1979
    32  = "21154"  # 1980
    33  = "16402"  # 1981
    34  = "13845"  # 1982
    35  = "14292"  # 1983
    36  = "14384"  # 1984
    37  = "13015"  # 1985
    38  = "13509"  # 1986
    39  = "12416"  # 1987
    40  = "12317"  # 1988
    41  = "12701"  # 1989
    42  = "9757"  # 1990
    43  = "6310.90684193707"  # 1991
    44  = "12965.0669266169"  # 1992
    45  = "19955.3143291455"  # 1993
    46  = "24532.4134922228"  # 1994
    47  = "27552.9609822888"  # 1995
    48  = "29612.6600989995"  # 1996
    49  = "31511.191303357"  # 1997
    50  = "33331.2986244481"  # 1998
    51  = "33351.1949905453"  # 1999
    52  = "35844.7475240851"  # 2000
    53  = "37225.8694104754"  # 2001
    54  = "39943.1481627285"  # 2002
    55  = "48809.7127623089"  # 2003
    56  = "55911.7935999109"  # 2004
    57  = "62959.6014217313"  # 2005
    58  = "68546.1428873113"  # 2006
    59  = "73096.8709734517"  # 2007
    60  = "75138.1196927745"  # 2008
    61  = "70072.729383949"  # 2009
    62  = "68865.3096214215"  # 2010
}
foreach ($row in $gdpUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $gdpUpdates[$row]
}

# Newly published years (2011-2016) appended to the bottom of the table.
$newRows = @(
    @{ Row = 63; Year = 2011; Value = "77126" }
    @{ Row = 64; Year = 2012; Value = "78801" }
    @{ Row = 65; Year = 2013; Value = "75256" }
    @{ Row = 66; Year = 2014; Value = "72508" }
    @{ Row = 67; Year = 2015; Value = "71354" }
    @{ Row = 68; Year = 2016; Value = "71010" }
)
foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = 414
    $ws.Cells.Item($r, 2).Value = "Kuwait"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $entry.Year
    $ws.Cells.Item($r, 5).Value = $entry.Value
}
